$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns are treated as text so that
# numeric-looking strings (e.g. "255.90", "1.00", "98.360.86") keep their
# exact literal formatting instead of being parsed as numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value2 = '98.360.86'
$ws.Cells.Item(2, 5).Value2 = '  -0.39%  '

$ws.Cells.Item(3, 4).Value2 = '3.420.50'
$ws.Cells.Item(3, 5).Value2 = '  +2.07%  '

$ws.Cells.Item(4, 5).Value2 = '  -0.03%  '

$ws.Cells.Item(5, 4).Value2 = '255.90'
$ws.Cells.Item(5, 5).Value2 = '  -1.59%  '

$ws.Cells.Item(6, 4).Value2 = '668.84'
$ws.Cells.Item(6, 5).Value2 = '  +2.22%  '

$ws.Cells.Item(7, 5).Value2 = '  -4.93%  '

$ws.Cells.Item(8, 4).Value2 = '0.439'
$ws.Cells.Item(8, 5).Value2 = '  -4.91%  '

$ws.Cells.Item(9, 5).Value2 = '  -2.41%  '

$ws.Cells.Item(10, 5).Value2 = '  +0.00%  '

$ws.Cells.Item(11, 4).Value2 = '3.415.93'
$ws.Cells.Item(11, 5).Value2 = '  +2.04%  '

$ws.Cells.Item(12, 5).Value2 = '  +3.13%  '

$ws.Cells.Item(13, 4).Value2 = '42.28'
$ws.Cells.Item(13, 5).Value2 = '  -2.64%  '

$ws.Cells.Item(14, 4).Value2 = '6.44'
$ws.Cells.Item(14, 5).Value2 = '  +14.65%  '

$ws.Cells.Item(15, 4).Value2 = '98.129.08'
$ws.Cells.Item(15, 5).Value2 = '  -0.14%  '

$ws.Cells.Item(16, 5).Value2 = '  -0.15%  '

$ws.Cells.Item(17, 4).Value2 = '4.049.63'
$ws.Cells.Item(17, 5).Value2 = '  +1.85%  '

$ws.Cells.Item(18, 4).Value2 = '9.01'
$ws.Cells.Item(18, 5).Value2 = '  +19.37%  '

$ws.Cells.Item(19, 4).Value2 = '3.422.45'
$ws.Cells.Item(19, 5).Value2 = '  +2.09%  '

$ws.Cells.Item(20, 4).Value2 = '0.577'
$ws.Cells.Item(20, 5).Value2 = '  +32.91%  '

$ws.Cells.Item(21, 4).Value2 = '17.75'
$ws.Cells.Item(21, 5).Value2 = '  +4.70%  '

$ws.Cells.Item(22, 4).Value2 = '11.10'
$ws.Cells.Item(22, 5).Value2 = '  +6.06%  '

$ws.Cells.Item(23, 5).Value2 = '  -4.29%  '

$ws.Cells.Item(24, 4).Value2 = '512.92'
$ws.Cells.Item(24, 5).Value2 = '  -4.16%  '

$ws.Cells.Item(25, 5).Value2 = '  -2.68%  '

$ws.Cells.Item(26, 4).Value2 = '6.71'
$ws.Cells.Item(26, 5).Value2 = '  +6.87%  '

$ws.Cells.Item(27, 4).Value2 = '101.86'
$ws.Cells.Item(27, 5).Value2 = '  +0.21%  '

$ws.Cells.Item(28, 4).Value2 = '12.90'
$ws.Cells.Item(28, 5).Value2 = '  +2.07%  '

$ws.Cells.Item(29, 4).Value2 = '3.606.84'
$ws.Cells.Item(29, 5).Value2 = '  +2.22%  '

$ws.Cells.Item(30, 5).Value2 = '  +1.17%  '

$ws.Cells.Item(31, 4).Value2 = '11.63'
$ws.Cells.Item(31, 5).Value2 = '  +5.29%  '

$ws.Cells.Item(32, 4).Value2 = '0.199'
$ws.Cells.Item(32, 5).Value2 = '  +3.22%  '

$ws.Cells.Item(33, 4).Value2 = '0.999'
$ws.Cells.Item(33, 5).Value2 = '  -0.16%  '

$ws.Cells.Item(34, 4).Value2 = '2.49'
$ws.Cells.Item(34, 5).Value2 = '  +18.41%  '

$ws.Cells.Item(35, 4).Value2 = '0.580'
$ws.Cells.Item(35, 5).Value2 = '  +7.41%  '

$ws.Cells.Item(36, 4).Value2 = '1.00'
$ws.Cells.Item(36, 5).Value2 = '  -0.73%  '

$ws.Cells.Item(37, 4).Value2 = '30.19'
$ws.Cells.Item(37, 5).Value2 = '  +2.49%  '

$ws.Cells.Item(38, 4).Value2 = '1.53'
$ws.Cells.Item(38, 5).Value2 = '  +15.23%  '

$ws.Cells.Item(39, 5).Value2 = '  +1.06%  '

$ws.Cells.Item(40, 4).Value2 = '541.73'
$ws.Cells.Item(40, 5).Value2 = '  +2.92%  '

$ws.Cells.Item(41, 5).Value2 = '  -1.97%  '

$ws.Cells.Item(42, 5).Value2 = '  -0.02%  '

$ws.Cells.Item(43, 4).Value2 = '0.881'
$ws.Cells.Item(43, 5).Value2 = '  +6.34%  '

$ws.Cells.Item(44, 4).Value2 = '24.71'
$ws.Cells.Item(44, 5).Value2 = '  +0.01%  '

$ws.Cells.Item(45, 2).Value2 = 'Cosmos'
$ws.Cells.Item(45, 3).Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(45, 4).Value2 = '9.05'
$ws.Cells.Item(45, 5).Value2 = '  +14.34%  '

$ws.Cells.Item(46, 2).Value2 = 'Filecoin'
$ws.Cells.Item(46, 3).Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(46, 4).Value2 = '5.88'
$ws.Cells.Item(46, 5).Value2 = '  +14.55%  '

$ws.Cells.Item(47, 2).Value2 = 'VeChain'
$ws.Cells.Item(47, 3).Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(47, 4).Value2 = '0.0434'
$ws.Cells.Item(47, 5).Value2 = '  -0.22%  '

$ws.Cells.Item(48, 2).Value2 = 'MantraDAO'
$ws.Cells.Item(48, 3).Value2 = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Cells.Item(48, 4).Value2 = '3.79'
$ws.Cells.Item(48, 5).Value2 = '  +0.64%  '

$ws.Cells.Item(49, 4).Value2 = '1.74'
$ws.Cells.Item(49, 5).Value2 = '  +15.41%  '

$ws.Cells.Item(50, 5).Value2 = '  -2.98%  '

$ws.Cells.Item(51, 4).Value2 = '54.11'
$ws.Cells.Item(51, 5).Value2 = '  +9.92%  '

# Restore the default cell style so no stray number-format styling is
# left behind on the price column (matches original unstyled cells).
$priceRange.Style = "Normal"
